# Realestate Update resale numbers 2024-01-06 23:46
# Append a new data row (row 28) to the CityResaleNum sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 28

# Columns A-D hold text values (date/time/weekday/week-code) that look
# numeric/date-like ("2024-01-06", "00"); format them as Text first so
# Excel stores them as literal strings instead of re-interpreting them
# as a date serial / number.
$ws.Range("A" + $row + ":D" + $row).NumberFormat = "@"

$ws.Cells.Item($row, 1).Value = "2024-01-06"
$ws.Cells.Item($row, 2).Value = "23:46:43"
$ws.Cells.Item($row, 3).Value = "Saturday"
$ws.Cells.Item($row, 4).Value = "00"

$ws.Cells.Item($row, 5).Value = 140547
$ws.Cells.Item($row, 6).Value = 143026
$ws.Cells.Item($row, 7).Value = 172337
$ws.Cells.Item($row, 8).Value = 147268
$ws.Cells.Item($row, 9).Value = -1
$ws.Cells.Item($row, 10).Value = 118395
$ws.Cells.Item($row, 11).Value = 224628
$ws.Cells.Item($row, 12).Value = 249367
$ws.Cells.Item($row, 13).Value = 185182
$ws.Cells.Item($row, 14).Value = 110416
$ws.Cells.Item($row, 15).Value = 40638
$ws.Cells.Item($row, 16).Value = 30810
$ws.Cells.Item($row, 17).Value = 72517
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 42276
$ws.Cells.Item($row, 20).Value = -1
